$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) cell format, used to restore
# style after forcing a text NumberFormat on numeric-looking price strings.
$plainStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = '26.948.36'
$ws.Range("D3").Value = '1.718.19'
$ws.Range("E3").Value = '  -2.67%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.45'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  -5.94%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4847'
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = '  +6.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3498'
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.09'
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07250'
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.045'
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  -4.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.89'
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = '  -3.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.855'
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").Value = '1.717.76'
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.853'
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = '  -4.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.57'
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = '  -6.28%  '
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.639'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("D23").Value = '27.000.83'
$ws.Range("E23").Value = '  -3.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.79'
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = '  -3.83%  '
$ws.Range("E25").Value = '  -3.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.68'
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = '  -5.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.90'
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = '  -1.24%  '
$ws.Range("D28").Value = '1.915.26'
$ws.Range("E28").Value = '  -2.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.068'
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = '  -4.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.91'
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = '  -2.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.026'
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = '  -4.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09284'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.586'
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.351'
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05879'
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = '  -3.99%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02177'
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = '  -4.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.447'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  +4.76%  '
$ws.Range("E38").Value = '  -7.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.1991'
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = '  -4.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6012'
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = '  -3.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.731'
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  -4.11%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.097'
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  -7.24%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.507'
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = '  -4.01%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.78'
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = '  -2.51%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.574'
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = '  -4.29%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5630'
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '117.98'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  -3.88%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.832'
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = '  -5.29%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.108'
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = '  -1.65%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06651'
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  -2.31%  '
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.000'
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = '  -0.05%  '
